$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 149.656361
$ws.Range("H2").Value = 448.969083
$ws.Range("I2").Value = 0.5921360794347563
$ws.Range("J2").Value = 0.5921360794347564
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.904209333333332
$ws.Range("N2").Value = 26.712628
$ws.Range("O2").Value = 0.3928890865119899
$ws.Range("P2").Value = 0.3928890865119899
$ws.Range("Q2").Value = 1332.571566408902
$ws.Range("R2").Value = 11993.14409768012
$ws.Range("S2").Value = 0.2326438033399125
$ws.Range("T2").Value = 0.2326438033399126

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 149.656361
$ws.Range("H3").Value = 448.969083
$ws.Range("I3").Value = 0.5921360794347563
$ws.Range("J3").Value = 0.5921360794347564
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.112038
$ws.Range("N3").Value = 27.336114
$ws.Range("O3").Value = 0.4020593128556135
$ws.Range("P3").Value = 0.4020593128556135
$ws.Range("Q3").Value = 1363.674448373718
$ws.Range("R3").Value = 12273.07003536346
$ws.Range("S3").Value = 0.2380738252145551
$ws.Range("T3").Value = 0.2380738252145552

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 149.656361
$ws.Range("H4").Value = 448.969083
$ws.Range("I4").Value = 0.5921360794347563
$ws.Range("J4").Value = 0.5921360794347564
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.64717
$ws.Range("N4").Value = 13.94151
$ws.Range("O4").Value = 0.2050516006323966
$ws.Range("P4").Value = 0.2050516006323966
$ws.Range("Q4").Value = 695.4785511483699
$ws.Range("R4").Value = 6259.30696033533
$ws.Range("S4").Value = 0.1214184508802887
$ws.Range("T4").Value = 0.1214184508802887

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 52.73412466666667
$ws.Range("H5").Value = 158.202374
$ws.Range("I5").Value = 0.208649853730866
$ws.Range("J5").Value = 0.208649853730866
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.904209333333332
$ws.Range("N5").Value = 26.712628
$ws.Range("O5").Value = 0.3928890865119899
$ws.Range("P5").Value = 0.3928890865119899
$ws.Range("Q5").Value = 469.5556850420969
$ws.Range("R5").Value = 4226.001165378872
$ws.Range("S5").Value = 0.08197625043318026
$ws.Range("T5").Value = 0.08197625043318027

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 52.73412466666667
$ws.Range("H6").Value = 158.202374
$ws.Range("I6").Value = 0.208649853730866
$ws.Range("J6").Value = 0.208649853730866
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.112038
$ws.Range("N6").Value = 27.336114
$ws.Range("O6").Value = 0.4020593128556135
$ws.Range("P6").Value = 0.4020593128556135
$ws.Range("Q6").Value = 480.515347859404
$ws.Range("R6").Value = 4324.638130734636
$ws.Range("S6").Value = 0.08388961681845626
$ws.Range("T6").Value = 0.08388961681845626

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 52.73412466666667
$ws.Range("H7").Value = 158.202374
$ws.Range("I7").Value = 0.208649853730866
$ws.Range("J7").Value = 0.208649853730866
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.64717
$ws.Range("N7").Value = 13.94151
$ws.Range("O7").Value = 0.2050516006323966
$ws.Range("P7").Value = 0.2050516006323966
$ws.Range("Q7").Value = 245.0644421271933
$ws.Range("R7").Value = 2205.57997914474
$ws.Range("S7").Value = 0.04278398647922949
$ws.Range("T7").Value = 0.0427839864792295

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 50.34932566666667
$ws.Range("H8").Value = 151.047977
$ws.Range("I8").Value = 0.1992140668343777
$ws.Range("J8").Value = 0.1992140668343777
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.904209333333332
$ws.Range("N8").Value = 26.712628
$ws.Range("O8").Value = 0.3928890865119899
$ws.Range("P8").Value = 0.3928890865119899
$ws.Range("Q8").Value = 448.3209355281728
$ws.Range("R8").Value = 4034.888419753556
$ws.Range("S8").Value = 0.07826903273889715
$ws.Range("T8").Value = 0.07826903273889717

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 50.34932566666667
$ws.Range("H9").Value = 151.047977
$ws.Range("I9").Value = 0.1992140668343777
$ws.Range("J9").Value = 0.1992140668343777
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.112038
$ws.Range("N9").Value = 27.336114
$ws.Range("O9").Value = 0.4020593128556135
$ws.Range("P9").Value = 0.4020593128556135
$ws.Range("Q9").Value = 458.784968749042
$ws.Range("R9").Value = 4129.064718741379
$ws.Range("S9").Value = 0.08009587082260215
$ws.Range("T9").Value = 0.08009587082260217

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 50.34932566666667
$ws.Range("H10").Value = 151.047977
$ws.Range("I10").Value = 0.1992140668343777
$ws.Range("J10").Value = 0.1992140668343777
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.64717
$ws.Range("N10").Value = 13.94151
$ws.Range("O10").Value = 0.2050516006323966
$ws.Range("P10").Value = 0.2050516006323966
$ws.Range("Q10").Value = 233.9818757583633
$ws.Range("R10").Value = 2105.83688182527
$ws.Range("S10").Value = 0.04084916327287837
$ws.Range("T10").Value = 0.04084916327287837
